$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.494.82"
$ws.Range("E2").Value = "  +3.92%  "
$ws.Range("D3").Value = "1.913.03"
$ws.Range("E3").Value = "  +2.30%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "333.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.85%  "
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4677"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.43%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4104"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.29%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.02"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08045"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.43%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.014"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.92%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.39"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.07%  "
$ws.Range("D13").Value = "1.949.33"
$ws.Range("E13").Value = "  +4.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.977"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.10%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.191"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.85%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "89.94"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.08%  "
$ws.Range("E18").Value = "  +1.57%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06585"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.74%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.14%  "
$ws.Range("D22").Value = "29.458.29"
$ws.Range("E22").Value = "  +3.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.574"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.29%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.51"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.208"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.78%  "
$ws.Range("D26").Value = "2.145.81"
$ws.Range("E26").Value = "  +2.68%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "155.47"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.88"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.765"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +9.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.144"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "117.38"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("E32").Value = "  +11.58%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09462"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.27%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.428"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.566"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.58%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.413"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06121"
$ws.Range("D37").Style = "Normal"
$ws.Range("E38").Value = "  +3.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.414"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.51%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.178"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5894"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.33%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1844"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.21"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.81%  "
$ws.Range("B44").Value = "WEMIXTOKEN"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.267"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.91%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.353"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.26%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.07504"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5573"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "12.17"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.932"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "113.36"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.2978"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +10.65%  "
